$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.677.32"
$ws.Range("E2").Value = "  -1.87%  "
$ws.Range("D3").Value = "1.740.62"
$ws.Range("E3").Value = "  -2.46%  "
$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = "  +0.69%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.88"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  -1.39%  "
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  +0.47%  "
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3872"
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = "  +1.09%  "
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3347"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  -2.64%  "
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.37"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  -4.83%  "
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.093"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  -5.82%  "
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07108"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  -4.39%  "
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.004"
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = "  +0.79%  "
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.67"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  -6.57%  "
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.047"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  -5.82%  "
$ws.Range("D15").Value = "1.741.00"
$ws.Range("E15").Value = "  -2.46%  "
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.914"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  -3.04%  "
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001041"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  -3.88%  "
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06586"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  -1.00%  "
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  +0.63%  "
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.43"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  -5.32%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.58"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -5.28%  "
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.141"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  -4.61%  "
$ws.Range("D23").Value = "27.711.61"
$ws.Range("E23").Value = "  -1.76%  "
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.42"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  -5.62%  "
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.399"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  +0.52%  "
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.54"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  -0.22%  "
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.59"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  -6.25%  "
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.261"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  -6.94%  "
$ws.Range("D29").Value = "1.939.17"
$ws.Range("E29").Value = "  -2.39%  "
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.260"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  -12.72%  "
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "126.72"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  -5.94%  "
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.052"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  +2.26%  "
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.710"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  -7.29%  "
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08670"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  -1.48%  "
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.88"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  -7.16%  "
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.065"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  -5.01%  "
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02234"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  -7.92%  "
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06016"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  -5.12%  "
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6355"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  -7.45%  "
$ws.Range("E40").Value = "  -1.73%  "
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2070"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  -5.23%  "
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  +0.54%  "
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.778"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  -6.76%  "
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.50"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  -5.40%  "
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.798"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  -1.38%  "
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5858"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  -7.41%  "
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.01"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  -5.51%  "
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.951"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  -6.96%  "
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06917"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  -7.12%  "
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.136"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  -4.18%  "
